# The Last Update 15-03-2024
# Refresh the Ligue 1 standings table (rows 2-19, columns A-G) with the
# latest team order and statistics, keeping every value stored as text
# (matching how the sheet already stores "1.3", "91%", etc.) instead of
# letting Excel auto-coerce numeric-looking strings into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    # Putting the cell back on the "Normal" style keeps the style index at
    # 0, exactly like the untouched cells, so we don't leave a stray
    # explicit format behind on the data cells.
    $r.Style = "Normal"
}

# # | Equipe | Cartões | Escanteios | 1.5+ | 2.5+ | Med. Gols
$data = @(
    @("1.",  "PSG",         "1.4", "5.4", "88%", "63%", "3.04"),
    @("2.",  "Brest",       "2.4", "4.4", "67%", "38%", "2.21"),
    @("3.",  "Monaco",      "2.4", "5.5", "83%", "67%", "3.25"),
    @("4.",  "LOSC Lille",  "1.8", "5.6", "76%", "37%", "2.25"),
    @("5.",  "Nice",        "1.8", "5.6", "46%", "21%", "1.67"),
    @("6.",  "Lens",        "2.4", "5.2", "72%", "42%", "2.38"),
    @("7.",  "Marseille",   "1.7", "5.8", "75%", "46%", "2.67"),
    @("8.",  "Rennes",      "1.8", "4.5", "78%", "52%", "2.71"),
    @("9.",  "Reims",       "1.8", "5.4", "79%", "54%", "2.58"),
    @("10.", "Toulouse",    "2.5", "4.3", "83%", "54%", "2.46"),
    @("11.", "Lyon",        "1.7", "4.5", "66%", "54%", "2.63"),
    @("12.", "Strasbourg",  "1.8", "3.4", "88%", "58%", "2.58"),
    @("13.", "Lorient",     "1.9", "3.4", "79%", "67%", "3.25"),
    @("14.", "Nantes",      "1.8", "5.1", "66%", "43%", "2.46"),
    @("15.", "Le Havre",    "1.8", "3.7", "67%", "54%", "2.29"),
    @("16.", "Montpellier", "2.1", "4.2", "76%", "54%", "2.50"),
    @("17.", "Metz",        "1.5", "3.9", "61%", "46%", "2.42"),
    @("18.", "Clermont",    "1.7", "4.3", "67%", "50%", "2.42")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]

    Set-TextValue "A$rowNum" $row[0]
    $ws.Range("B$rowNum").Value = $row[1]
    Set-TextValue "C$rowNum" $row[2]
    Set-TextValue "D$rowNum" $row[3]
    Set-TextValue "E$rowNum" $row[4]
    Set-TextValue "F$rowNum" $row[5]
    Set-TextValue "G$rowNum" $row[6]
}
